$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add four new log rows (6-9), cloning row 5's formatting so the new
# rows pick up the same date / text styles already used in the table.
$ws.Rows(5).Copy()
$ws.Rows(6).Insert()
$ws.Rows(5).Copy()
$ws.Rows(7).Insert()
$ws.Rows(5).Copy()
$ws.Rows(8).Insert()
$ws.Rows(5).Copy()
$ws.Rows(9).Insert()

# Row 6 - 04 / JEB / System Design document updated
$ws.Cells.Item(6, 1).Value = 41710
$ws.Cells.Item(6, 2).Value = "04"
$ws.Cells.Item(6, 3).Value = "JEB"
$ws.Cells.Item(6, 4).Value = "System Design document updated"
$ws.Cells.Item(6, 5).Value = "Done"

# Row 7 - 05 / JEB / Test Cases Document for Scheduler Updated
$ws.Cells.Item(7, 1).Value = 41710
$ws.Cells.Item(7, 2).Value = "05"
$ws.Cells.Item(7, 3).Value = "JEB"
$ws.Cells.Item(7, 4).Value = "Test Cases Document for Scheduler Updated"
$ws.Cells.Item(7, 5).Value = "Done"

# Row 8 - 06 / JEB / Traceability document for Scheduler Updated
$ws.Cells.Item(8, 1).Value = 41710
$ws.Cells.Item(8, 2).Value = "06"
$ws.Cells.Item(8, 3).Value = "JEB"
$ws.Cells.Item(8, 4).Value = "Traceability document for Scheduler Updated"
$ws.Cells.Item(8, 5).Value = "Done"

# Row 9 - 07 / JEB / Report Documents updated
$ws.Cells.Item(9, 1).Value = 41711
$ws.Cells.Item(9, 2).Value = "07"
$ws.Cells.Item(9, 3).Value = "JEB"
$ws.Cells.Item(9, 4).Value = "Report Documents updated"
$ws.Cells.Item(9, 5).Value = "Done"

# Keep the sheet's dimension/selection in sync with the new last row.
$ws.Range("E10").Select() | Out-Null
